$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,3,4) got cyclically rotated:
#   new row2 = old row4
#   new row3 = old row2
#   new row4 = old row3
# (for the columns that actually hold per-record data)
$cols = @("A","B","D","E","F","G","H","Q","R")

$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($c in $cols) {
    $row2[$c] = $ws.Range($c + "2").Value2
    $row3[$c] = $ws.Range($c + "3").Value2
    $row4[$c] = $ws.Range($c + "4").Value2
}

foreach ($c in $cols) {
    $ws.Range($c + "2").Value2 = $row4[$c]
    $ws.Range($c + "3").Value2 = $row2[$c]
    $ws.Range($c + "4").Value2 = $row3[$c]
}
